# Add a 't+3' column (column C) of predicted values to the worksheet,
# mirroring the existing 't+1'/'t+2' style columns A and B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Predicted values for the new 't+3' column, row by row (row 1 is the header).
$tplus3 = @(
    2,
    -0.3093622831163451,
    0.2070764539413797,
    0.1431646595838115,
    -0.0758442733719054,
    0.1446280657153557,
    -0.4512205850314532,
    -0.2584775454032165,
    -0.4105889929375005,
    0.3588729555905856,
    -0.2385017336592683,
    -0.08619714041918237,
    -0.0118996560469052,
    0.05218774876008093,
    -0.01411693445400752,
    0.3464168503908811,
    0.5981719936045718,
    0.03604968763938152,
    0.3972573364645569,
    0.2935978017111165,
    0.519205053466896,
    0.3776076824557239,
    -0.06027967599379161,
    4.564317052094395,
    0.4732361303223824,
    0.3793329866124788,
    0.2933557710250015,
    1.069421889075762,
    5.420665640241038,
    0.9469909350344954,
    -0.2215551915663329,
    0.7849417705186705,
    0.9049151329377539,
    -0.665537397882463,
    0.8139125278147399,
    0.7679010701018232,
    0.7405970643419434,
    0.7172636952775983,
    0.5813612769236129,
    0.7530707326558662,
    0.5544038253167535,
    0.6903968288805641,
    0.7086778126739838,
    0.6742917402476475,
    0.6582783396307987,
    -1.262984666059046,
    -0.9787528595606133,
    -0.8725351967226652,
    -0.6414807974495064,
    -0.04809784013177344,
    -0.8643473716963861,
    -0.8643473716963861,
    -1.096695474552693,
    -0.1862255579559371,
    -0.9981244871884227,
    -0.8834438227905759,
    -0.9329141727921776,
    -1.126906824887453,
    -0.8468034672072153,
    -0.4960781952889813,
    0.3717217949089451,
    -1.206867720350009,
    -0.7371004760351398,
    -0.9007103462505319,
    -0.109955290616581,
    -0.7887724666920192,
    -0.7605940848606052
)

# Copy the header cell's formatting (bold, bordered, centered) into C1
# before writing values, so the new column matches columns A and B.
$ws.Range("B1").Copy($ws.Range("C1")) | Out-Null

# Loop over every row and fill in the new column C with the 't+3' values.
for ($i = 0; $i -lt $tplus3.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $tplus3[$i]
}

